$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12: fill in previously empty pareto extrapolation cells
$ws.Range("I12").Value = 0.3222525556017057
$ws.Range("J12").Value = 0.05337495275340508
$ws.Range("K12").Value = -0.4103541666894833
$ws.Range("L12").Value = 1.728572074161355

# Row 17: fill in previously empty pareto extrapolation cells
$ws.Range("I17").Value = 0.4887390112796265
$ws.Range("J17").Value = 0.1186085878524447
$ws.Range("K17").Value = 0.1721465981455686
$ws.Range("L17").Value = 2.224517934906989
